$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new date columns (JV1:JZ1), same style/number-format as JU1 ---
# Copy JU1's formatting (date style) across the new cells first, then overwrite
# the values, so the new cells pick up the existing date-style index instead
# of Excel minting a brand-new (duplicate) style entry.
$ws.Range("JU1").Copy($ws.Range("JV1:JZ1"))

$ws.Range("JV1").Value = 44118
$ws.Range("JW1").Value = 44119
$ws.Range("JX1").Value = 44120
$ws.Range("JY1").Value = 44121
$ws.Range("JZ1").Value = 44122

# --- Row 2 ---
$ws.Range("JV2").Value = 52.71
$ws.Range("JW2").Value = 54.51
$ws.Range("JX2").Value = 51.44
$ws.Range("JY2").Value = 65.09
$ws.Range("JZ2").Value = 50.49

# --- Row 3 ---
$ws.Range("JV3").Value = 37.35
$ws.Range("JW3").Value = 37.83
$ws.Range("JX3").Value = 42.93
$ws.Range("JY3").Value = 39.27
$ws.Range("JZ3").Value = 35.45

# --- Row 4 ---
$ws.Range("JV4").Value = 60.84
$ws.Range("JW4").Value = 65.83
$ws.Range("JX4").Value = 61.49
$ws.Range("JY4").Value = 66.89
$ws.Range("JZ4").Value = 52.7

# --- Row 5 ---
$ws.Range("JV5").Value = 63.31
$ws.Range("JW5").Value = 61.75
$ws.Range("JX5").Value = 65.57
$ws.Range("JY5").Value = 63.15
$ws.Range("JZ5").Value = 52

# --- View state: scroll back to the top-left and select the full used rows
#     (mirrors the author re-selecting rows 1:5 after pasting the new data) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Rows("1:5").Select() | Out-Null
